# Variable_Index.xlsx edit:
#  1. Insert a new worksheet "CCES_Ver1 Variables" between "List of Variables"
#     and "Section Questions", populated with a snapshot of the variable list.
#  2. Rename the shared string used by 3 Year Overview!H40 from
#     "Prohibit > 20 Weeks_16" to "Prohibit MoreThan20 Weeks_16".
#  3. Make "3 Year Overview" the active tab / selected sheet, and nudge a few
#     sheet selections to match the saved UI state.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("3 Year Overview")
$ws2 = $wb.Worksheets.Item("List of Variables")
$ws3 = $wb.Worksheets.Item("Section Questions")

# ---------------------------------------------------------------------------
# 1. New sheet "CCES_Ver1 Variables", inserted right before "Section Questions"
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($ws3)
$newSheet.Name = "CCES_Ver1 Variables"

# Header row
$newSheet.Range("A2").Value = "Class"
$newSheet.Range("A2").Font.Bold = $true
$newSheet.Range("B2").Value = "SQL Index"
$newSheet.Range("B2").Font.Bold = $true

# Data rows: (Class label, SQL Index variable name)
$rows = @(
    @(3,  "ID",            "V101_16"),
    @(4,  "Location",       "Zipcode_16"),
    @(5,  "Location",       "State_16"),
    @(6,  "Location",       "CountyFips_16"),
    @(7,  "Location",       "CountyName_16"),
    @(8,  "Self",            "Birth Year_16"),
    @(9,  "Self",            "Gender_16"),
    @(10, "Self",            "Education_16"),
    @(11, "Self",            "Marrital Status_16"),
    @(12, "Self Family",     "Children < 18 yrs_16"),
    @(14, "Crime",           "Gun Background Checks_16"),
    @(15, "Crime",           "Prohibit Publication_16"),
    @(16, "Crime",           "Ban Assult Weapons_16"),
    @(17, "Crime",           "Make CCP Easier_16"),
    @(18, "Abortion",        "Always Allow Choice_16"),
    @(19, "Abortion",        "Rape, Incest, or Health_16"),
    @(20, "Abortion",        "Prohibit > 20 Weeks_16"),
    @(21, "Abortion",        "Employers decline benefits_16"),
    @(22, "Abortion",        "Prohibit Fed Funds_16"),
    @(23, "Abortion",        "Illegal in all circumstances_16"),
    @(24, "Gay Marriage",    "Gay Marriage_16")
)

foreach ($row in $rows) {
    $r = $row[0]
    $newSheet.Range("A$r").Value = $row[1]
    $newSheet.Range("B$r").Value = $row[2]
}

$newSheet.Columns("A").ColumnWidth = 12.5546875
$newSheet.Columns("B").ColumnWidth = 37.44140625

$newSheet.Range("C3").Select()

# ---------------------------------------------------------------------------
# 2. Rename the shared string referenced by H40 on "3 Year Overview"
# ---------------------------------------------------------------------------
$ws1.Range("H40").Value = "Prohibit MoreThan20 Weeks_16"

# ---------------------------------------------------------------------------
# 3. Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------
$ws2.Range("A22:XFD22").Select()
$ws3.Range("K6").Select()

$ws1.Range("H42").Select()
$ws1.Activate()

Write-Output "edit complete"
